$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M6").Value = -188
$ws.Range("I6").Value = 100
$ws.Range("K6").Value = 300
$ws.Range("H6").Value = 100
$ws.Range("M33").Value = -205.0909
$ws.Range("I33").Value = 434.0909
$ws.Range("K33").Value = 434.0909
$ws.Range("H33").Value = 451.30768
$ws.Range("H40").Value = 2209.75
$ws.Range("I40").Value = 3300
$ws.Range("K40").Value = 3300
$ws.Range("M40").Value = -3125
$ws.Range("H41").Value = 239.5
$ws.Range("K41").Value = 92.5
$ws.Range("M41").Value = 347.5
$ws.Range("I41").Value = 92.5
$ws.Range("I61").Value = 583
$ws.Range("H61").Value = 987.8570999999999
$ws.Range("M61").Value = -1577
$ws.Range("K61").Value = 1749
$ws.Range("K76").Value = 333337000
$ws.Range("I76").Value = 333337000
$ws.Range("M76").Value = -333336685
$ws.Range("H76").Value = 250004220
$ws.Range("M79").Value = -333335908
$ws.Range("H79").Value = 250004220
$ws.Range("I79").Value = 333337000
$ws.Range("K79").Value = 333337000
$ws.Range("M80").Value = -2087.2858
$ws.Range("N80").Value = -6112.6666
$ws.Range("H80").Value = 1221.8125
$ws.Range("I80").Value = 1028.4286
$ws.Range("L80").Value = 4116.6666
$ws.Range("K80").Value = 3085.2858
$ws.Range("J80").Value = 1372.2222
$ws.Range("K82").Value = 28500
$ws.Range("H82").Value = 9500
$ws.Range("M82").Value = -28094
$ws.Range("I82").Value = 9500
$ws.Range("N83").Value = -22333.9998
$ws.Range("I83").Value = 1028.4286
$ws.Range("K83").Value = 9255.857399999999
$ws.Range("L83").Value = 12349.9998
$ws.Range("M83").Value = -4263.857399999999
$ws.Range("J83").Value = 1372.2222
$ws.Range("H83").Value = 1221.8125
$ws.Range("I85").Value = 9500
$ws.Range("H85").Value = 9500
$ws.Range("M85").Value = -27096
$ws.Range("K85").Value = 28500
$ws.Range("N111").Value = -27882.5
$ws.Range("H111").Value = 4331.75
$ws.Range("J111").Value = 7249.5
$ws.Range("L111").Value = 21748.5
$ws.Range("L112").Value = 6934.799999999999
$ws.Range("N112").Value = -9150.799999999999
$ws.Range("H112").Value = 2526.348
$ws.Range("J112").Value = 2311.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -1361.32
$ws.Range("K2").Value = 1474.32
$ws.Range("I2").Value = 1474.32
$ws.Range("H2").Value = 2805.4138
$ws.Range("J2").Value = 11124.75
$ws.Range("L2").Value = 11124.75
$ws.Range("N2").Value = -11350.75
$ws.Range("H32").Value = 1711566.8
$ws.Range("I32").Value = 1749.6394
$ws.Range("M32").Value = -1462.6394
$ws.Range("K32").Value = 1749.6394
$ws.Range("I61").Value = 8185.5
$ws.Range("H61").Value = 5300.979
$ws.Range("M61").Value = -7973.5
$ws.Range("K61").Value = 8185.5
$ws.Range("K116").Value = 1474.32
$ws.Range("M116").Value = 819.6800000000001
$ws.Range("H116").Value = 2805.4138
$ws.Range("J116").Value = 11124.75
$ws.Range("N116").Value = -15712.75
$ws.Range("I116").Value = 1474.32
$ws.Range("L116").Value = 11124.75
$ws.Range("I132").Value = 752092.4
$ws.Range("M132").Value = -2253747.2
$ws.Range("K132").Value = 2256277.2
$ws.Range("H132").Value = 672395
$ws.Range("I136").Value = 8185.5
$ws.Range("M136").Value = -22006.5
$ws.Range("H136").Value = 5300.979
$ws.Range("K136").Value = 24556.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 1474.32
$ws.Range("M3").Value = -1360.32
$ws.Range("H3").Value = 2805.4138
$ws.Range("K3").Value = 1474.32
$ws.Range("L3").Value = 11124.75
$ws.Range("J3").Value = 11124.75
$ws.Range("N3").Value = -11352.75
$ws.Range("M105").Value = -613.5264000000002
$ws.Range("K105").Value = 2360.5264
$ws.Range("H105").Value = 2342.5
$ws.Range("I105").Value = 2360.5264
$ws.Range("M134").Value = -5013551.4
$ws.Range("H134").Value = 1321942
$ws.Range("K134").Value = 5016086.4
$ws.Range("I134").Value = 1672028.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J94").Value = 5781.375
$ws.Range("L94").Value = 5781.375
$ws.Range("K94").Value = 100001740
$ws.Range("I94").Value = 100001740
$ws.Range("H94").Value = 55559092
$ws.Range("M94").Value = -100001289
$ws.Range("N94").Value = -6683.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J7").Value = 0
$ws.Range("I7").Value = 10545.5
$ws.Range("M7").Value = -31524.5
$ws.Range("K7").Value = 31636.5
$ws.Range("H7").Value = 10545.5
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H32").Value = 846.8570999999999
$ws.Range("J32").Value = 1300.6666
$ws.Range("N32").Value = -4467.9998
$ws.Range("L32").Value = 3901.9998
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("K122").Value = 6171.428699999999
$ws.Range("I122").Value = 685.7143
$ws.Range("H122").Value = 4627.96
$ws.Range("M122").Value = -3721.428699999999
$ws.Range("L137").Value = 10823.625
$ws.Range("J137").Value = 3607.875
$ws.Range("N137").Value = -21023.625
$ws.Range("M137").Value = 821.5715999999993
$ws.Range("K137").Value = 4278.428400000001
$ws.Range("I137").Value = 1426.1428
$ws.Range("H137").Value = 2589.7334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I102").Value = 2255080.8
$ws.Range("K102").Value = 2255080.8
$ws.Range("M102").Value = -2253458.8
$ws.Range("H102").Value = 1190393.1
$ws.Range("H126").Value = 20842890
$ws.Range("L126").Value = 45135.273
$ws.Range("N126").Value = -50075.273
$ws.Range("J126").Value = 15045.091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I61").Value = 5498.2104
$ws.Range("H61").Value = 6506.5713
$ws.Range("M61").Value = -5296.2104
$ws.Range("L61").Value = 8635.333000000001
$ws.Range("J61").Value = 8635.333000000001
$ws.Range("K61").Value = 5498.2104
$ws.Range("N61").Value = -9039.333000000001
$ws.Range("K82").Value = 2209.7
$ws.Range("H82").Value = 3804.2727
$ws.Range("M82").Value = -1848.7
$ws.Range("I82").Value = 2209.7
$ws.Range("I85").Value = 2209.7
$ws.Range("H85").Value = 3804.2727
$ws.Range("M85").Value = -961.6999999999998
$ws.Range("K85").Value = 2209.7
$ws.Range("H93").Value = 1385.2307
$ws.Range("M93").Value = -197.0667000000001
$ws.Range("L93").Value = 1185.7778
$ws.Range("N93").Value = -3681.7778
$ws.Range("I93").Value = 1445.0667
$ws.Range("K93").Value = 1445.0667
$ws.Range("J93").Value = 1185.7778
$ws.Range("L113").Value = 8635.333000000001
$ws.Range("J113").Value = 8635.333000000001
$ws.Range("K113").Value = 5498.2104
$ws.Range("H113").Value = 6506.5713
$ws.Range("N113").Value = -12975.333
$ws.Range("M113").Value = -3328.2104
$ws.Range("I113").Value = 5498.2104
$ws.Range("I132").Value = 5296
$ws.Range("M132").Value = -13358
$ws.Range("K132").Value = 15888
$ws.Range("H132").Value = 6206

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15560.2
$ws.Range("M62").Value = -13464
$ws.Range("I62").Value = 14088
$ws.Range("K62").Value = 14088
$ws.Range("N62").Value = -17789.666
$ws.Range("L62").Value = 16541.666
$ws.Range("J62").Value = 16541.666
$ws.Range("K65").Value = 70440
$ws.Range("M65").Value = -67320
$ws.Range("J65").Value = 16541.666
$ws.Range("L65").Value = 82708.33
$ws.Range("I65").Value = 14088
$ws.Range("N65").Value = -88948.33
$ws.Range("H65").Value = 15560.2
$ws.Range("K107").Value = 8965.636200000001
$ws.Range("M107").Value = -7045.636200000001
$ws.Range("I107").Value = 2988.5454
$ws.Range("H107").Value = 3232.1765
$ws.Range("L113").Value = 3409.5
$ws.Range("J113").Value = 1136.5
$ws.Range("H113").Value = 5377421.5
$ws.Range("N113").Value = -7749.5
$ws.Range("I132").Value = 3626.1924
$ws.Range("M132").Value = -8348.5772
$ws.Range("K132").Value = 10878.5772
$ws.Range("H132").Value = 4778.926
$ws.Range("N132").Value = -109310
$ws.Range("J132").Value = 34750
$ws.Range("L132").Value = 104250
